# Update the cached "datetimeFigureOut" Date placeholder text from
# 8/20/2020 to 10/27/2020 on the Slide Master and every Slide Layout
# (ppPlaceholderDate = 16 in the PpPlaceholderType enum).

$p = $ppt.ActivePresentation

$oldDate = "8/20/2020"
$newDate = "10/27/2020"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePh = $false
            if ($shp.Type -eq 14) {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePh = $true
                }
            }
            if ($isDatePh) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Walk every Slide Master (Design) in the deck, plus every Slide Layout
# (CustomLayout) hanging off each master.
for ($di = 1; $di -le $p.Designs.Count; $di++) {
    $master = $p.Designs.Item($di).SlideMaster

    Update-DatePlaceholder $master.Shapes

    for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
        $layout = $master.CustomLayouts.Item($li)
        Update-DatePlaceholder $layout.Shapes
    }
}
